# Updates cryptos list values (Price and Volume(1h) columns) to latest
# snapshot, matching commit "Updated cryptos list on Tue Jun 11 04:30:08 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.977.29"
$ws.Range("E2").Value = "  -2.34%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.553.81"
$ws.Range("E3").Value = "  -3.57%  "

# Row 4
$ws.Range("E4").Value = "  +0.11%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "622.66"
$ws.Range("E5").Value = "  -6.58%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "155.17"
$ws.Range("E6").Value = "  -2.93%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.546.79"
$ws.Range("E7").Value = "  -3.70%  "

# Row 8
$ws.Range("E8").Value = "  +0.28%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.487"
$ws.Range("E9").Value = "  -2.53%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.141"
$ws.Range("E10").Value = "  -3.14%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.95"
$ws.Range("E11").Value = "  -2.61%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.433"
$ws.Range("E12").Value = "  -1.97%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000225"
$ws.Range("E13").Value = "  -3.49%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.160.69"
$ws.Range("E14").Value = "  -3.37%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "32.16"
$ws.Range("E15").Value = "  -2.33%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.551.10"
$ws.Range("E16").Value = "  -3.74%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "68.050.01"
$ws.Range("E17").Value = "  -2.22%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.116"
$ws.Range("E18").Value = "  -0.73%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.43"
$ws.Range("E19").Value = "  -0.63%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.64"
$ws.Range("E20").Value = "  -3.36%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "460.70"
$ws.Range("E21").Value = "  -1.96%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.75"
$ws.Range("E22").Value = "  -0.05%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.643"
$ws.Range("E23").Value = "  -0.50%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "77.90"
$ws.Range("E24").Value = "  -2.26%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.704.09"
$ws.Range("E25").Value = "  -3.31%  "

# Row 26
$ws.Range("E26").Value = "  -0.01%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.69"
$ws.Range("E27").Value = "  -2.15%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0000116"
$ws.Range("E28").Value = "  -8.46%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.39"
$ws.Range("E29").Value = "  -7.32%  "

# Row 30
$ws.Range("E30").Value = "  -3.07%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.63"
$ws.Range("E31").Value = "  -4.22%  "

# Row 32
$ws.Range("E32").Value = "  +0.04%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "26.03"
$ws.Range("E33").Value = "  -2.64%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.91"
$ws.Range("E34").Value = "  -4.58%  "

# Row 35
$ws.Range("B35").Value = "Kaspa"
$ws.Range("C35").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.158"
$ws.Range("E35").Value = "  -4.61%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.19"
$ws.Range("E36").Value = "  -4.46%  "

# Row 37
$ws.Range("B37").Value = "RenzoRestakedETH"
$ws.Range("C37").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.557.61"
$ws.Range("E37").Value = "  -3.31%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.09"
$ws.Range("E38").Value = "  -4.31%  "

# Row 39
$ws.Range("E39").Value = "  +0.03%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "178.08"
$ws.Range("E40").Value = "  +0.75%  "

# Row 41
$ws.Range("E41").Value = "  +0.09%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.63"
$ws.Range("E42").Value = "  -7.59%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0885"
$ws.Range("E43").Value = "  -2.40%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.12"
$ws.Range("E44").Value = "  -5.61%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.894"
$ws.Range("E45").Value = "  -4.21%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "45.91"
$ws.Range("E46").Value = "  -2.40%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "28.49"
$ws.Range("E47").Value = "  +3.59%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.59"
$ws.Range("E48").Value = "  -5.64%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.73"
$ws.Range("E49").Value = "  -1.59%  "

# Row 50
$ws.Range("E50").Value = "  -5.27%  "

# Row 51
$ws.Range("E51").Value = "  -4.79%  "
